$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "EXN(2.4700356603552773, -0.1013161350166013, 1.890546768124959)"
$ws.Range("C2").Value = "NCT(2.6568777062030153, 1.59530627445154, 1.365815049779554, 4.250801857631765)"
$ws.Range("D2").Value = "NCT(2.5226423274199843, 1.2592767947879548, -0.36605784444340883, 2.637145461244753)"
$ws.Range("E2").Value = "MIE(6.8118061426931895, 3.987097137097183, -12.03678734558644, 17.67260118347503)"
